$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New example rows to append beneath the existing data (rows 5-9).
$rows = @(
    @("bar {abc}Bar", "Lorem ...", $false, "asd", "dfg"),
    @("bar{abc} Bar", "Lorem ...", $false, "asd", "dfg"),
    @("{123}", "Lorem ...", $false, "asd", "dfg"),
    @("{123} {abc}", "Lorem ...", $false, "asd", "dfg"),
    @("foo {abc} bar", "Lorem ...", $false, "asd", "dfg")
)

$r = 5
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Widen column A slightly to fit the new, longer example values
# (this also splits the former A:B "col" run into separate A and B entries).
$ws.Range("A1:A9").ColumnWidth = 10.666666666666666

# Move / expand the selection to the newly added last row.
$ws.Range("A9:XFD9").Select()
